# Scheduled-runner refresh of the Goblin Profits pricing snapshot.
# Updates the market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) on a handful of rows across each job sheet to the latest
# pulled values. CRP row 57's LeveProfitHQ (N57) is cleared because that
# leve no longer has an HQ turn-in profit figure this cycle.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3300.1428
$ws.Range("I18").Value = 1725.25
$ws.Range("J18").Value = 5400
$ws.Range("K18").Value = 1725.25
$ws.Range("L18").Value = 5400
$ws.Range("M18").Value = -1441.25
$ws.Range("N18").Value = -5968

$ws.Range("H64").Value = 8983.040000000001
$ws.Range("I64").Value = 6429.3335
$ws.Range("J64").Value = 9789.474
$ws.Range("K64").Value = 6429.3335
$ws.Range("L64").Value = 9789.474
$ws.Range("M64").Value = -6181.3335
$ws.Range("N64").Value = -10285.474

$ws.Range("H67").Value = 8983.040000000001
$ws.Range("I67").Value = 6429.3335
$ws.Range("J67").Value = 9789.474
$ws.Range("K67").Value = 6429.3335
$ws.Range("L67").Value = 9789.474
$ws.Range("M67").Value = -5571.3335
$ws.Range("N67").Value = -11505.474

$ws.Range("H80").Value = 1686.8
$ws.Range("I80").Value = 747.625
$ws.Range("J80").Value = 2312.9167
$ws.Range("K80").Value = 2242.875
$ws.Range("L80").Value = 6938.750100000001
$ws.Range("M80").Value = -1244.875
$ws.Range("N80").Value = -8934.750100000001

$ws.Range("H83").Value = 1686.8
$ws.Range("I83").Value = 747.625
$ws.Range("J83").Value = 2312.9167
$ws.Range("K83").Value = 6728.625
$ws.Range("L83").Value = 20816.2503
$ws.Range("M83").Value = -1736.625
$ws.Range("N83").Value = -30800.2503

$ws.Range("H88").Value = 7776.4165
$ws.Range("I88").Value = 3663.4
$ws.Range("J88").Value = 10714.286
$ws.Range("K88").Value = 3663.4
$ws.Range("L88").Value = 10714.286
$ws.Range("M88").Value = -3257.4
$ws.Range("N88").Value = -11526.286

$ws.Range("H91").Value = 7776.4165
$ws.Range("I91").Value = 3663.4
$ws.Range("J91").Value = 10714.286
$ws.Range("K91").Value = 3663.4
$ws.Range("L91").Value = 10714.286
$ws.Range("M91").Value = -2259.4
$ws.Range("N91").Value = -13522.286

$ws.Range("H132").Value = 2011.6316
$ws.Range("I132").Value = 1762.1111
$ws.Range("K132").Value = 5286.3333
$ws.Range("M132").Value = -2756.3333

$ws.Range("H135").Value = 1102.2307
$ws.Range("I135").Value = 344.2
$ws.Range("K135").Value = 3097.8
$ws.Range("M135").Value = -562.7999999999997

$ws.Range("H137").Value = 2586.2593
$ws.Range("I137").Value = 2240.818
$ws.Range("J137").Value = 2823.75
$ws.Range("K137").Value = 6722.454000000001
$ws.Range("L137").Value = 8471.25
$ws.Range("M137").Value = -4172.454000000001
$ws.Range("N137").Value = -13571.25

$ws.Range("H141").Value = 6713.7334
$ws.Range("J141").Value = 17988.5
$ws.Range("L141").Value = 53965.5
$ws.Range("N141").Value = -64325.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 66345.06
$ws.Range("I32").Value = 70214.37
$ws.Range("J32").Value = 24335.428
$ws.Range("K32").Value = 70214.37
$ws.Range("L32").Value = 24335.428
$ws.Range("M32").Value = -69927.37
$ws.Range("N32").Value = -24909.428

$ws.Range("H92").Value = 33280
$ws.Range("J92").Value = 33280
$ws.Range("L92").Value = 33280
$ws.Range("N92").Value = -38272

$ws.Range("H110").Value = 1580.25
$ws.Range("I110").Value = 1580.25
$ws.Range("K110").Value = 1580.25
$ws.Range("M110").Value = 464.75

$ws.Range("H132").Value = 7055.091
$ws.Range("I132").Value = 7972.9414
$ws.Range("K132").Value = 23918.8242
$ws.Range("M132").Value = -21388.8242

$ws.Range("H133").Value = 77748.5
$ws.Range("J133").Value = 77748.5
$ws.Range("L133").Value = 77748.5
$ws.Range("N133").Value = -82808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4503.7026
$ws.Range("I20").Value = 5516.9614
$ws.Range("K20").Value = 5516.9614
$ws.Range("M20").Value = -5269.9614

$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 2006647.8
$ws.Range("J9").Value = 2006647.8
$ws.Range("L9").Value = 2006647.8
$ws.Range("N9").Value = -2006983.8

$ws.Range("H16").Value = 898.38464
$ws.Range("I16").Value = 761.8182
$ws.Range("K16").Value = 761.8182
$ws.Range("M16").Value = -474.8182

$ws.Range("H57").Value = 26000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H113").Value = 898.38464
$ws.Range("I113").Value = 761.8182
$ws.Range("K113").Value = 761.8182
$ws.Range("M113").Value = 1408.1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5432.185
$ws.Range("J131").Value = 8451.375
$ws.Range("L131").Value = 25354.125
$ws.Range("N131").Value = -35434.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 12000
$ws.Range("J44").Value = 12000
$ws.Range("L44").Value = 12000
$ws.Range("N44").Value = -13192

$ws.Range("H99").Value = 3436.5
$ws.Range("I99").Value = 718.6
$ws.Range("J99").Value = 7966.3335
$ws.Range("K99").Value = 718.6
$ws.Range("L99").Value = 7966.3335
$ws.Range("M99").Value = 1527.4
$ws.Range("N99").Value = -12458.3335

$ws.Range("H113").Value = 71444490
$ws.Range("I113").Value = 333367940
$ws.Range("J113").Value = 10818.182
$ws.Range("K113").Value = 333367940
$ws.Range("L113").Value = 10818.182
$ws.Range("M113").Value = -333365770
$ws.Range("N113").Value = -15158.182

$ws.Range("H132").Value = 2544.1
$ws.Range("I132").Value = 1943.3182
$ws.Range("K132").Value = 5829.9546
$ws.Range("M132").Value = -3299.9546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4618.143
$ws.Range("I61").Value = 1404.625
$ws.Range("K61").Value = 1404.625
$ws.Range("M61").Value = -1202.625

$ws.Range("H113").Value = 4618.143
$ws.Range("I113").Value = 1404.625
$ws.Range("K113").Value = 1404.625
$ws.Range("M113").Value = 765.375

$ws.Range("H132").Value = 3082000.2
$ws.Range("I132").Value = 2555.4443
$ws.Range("J132").Value = 10010751
$ws.Range("K132").Value = 7666.3329
$ws.Range("L132").Value = 30032253
$ws.Range("M132").Value = -5136.3329
$ws.Range("N132").Value = -30037313

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3975
$ws.Range("I81").Value = 3966.6667
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 7933.3334
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -6872.3334
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 3975
$ws.Range("I84").Value = 3966.6667
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 39666.667
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -34362.667
$ws.Range("N84").Value = -50608

$ws.Range("H96").Value = 6166.5
$ws.Range("I96").Value = 5999.5
$ws.Range("K96").Value = 5999.5
$ws.Range("M96").Value = -4626.5

$ws.Range("H100").Value = 744.8095
$ws.Range("I100").Value = 533.9286
$ws.Range("K100").Value = 1067.8572
$ws.Range("M100").Value = -526.8571999999999

$ws.Range("H107").Value = 1028.2778
$ws.Range("J107").Value = 1249.5
$ws.Range("L107").Value = 3748.5
$ws.Range("N107").Value = -7588.5

$ws.Range("H113").Value = 1595.091
$ws.Range("I113").Value = 1499.75
$ws.Range("J113").Value = 1649.5714
$ws.Range("K113").Value = 4499.25
$ws.Range("L113").Value = 4948.7142
$ws.Range("M113").Value = -2329.25
$ws.Range("N113").Value = -9288.7142

$ws.Range("H132").Value = 9247.509
$ws.Range("I132").Value = 6719.7334
$ws.Range("J132").Value = 20622.5
$ws.Range("K132").Value = 20159.2002
$ws.Range("L132").Value = 61867.5
$ws.Range("M132").Value = -17629.2002
$ws.Range("N132").Value = -66927.5

$ws.Range("H133").Value = 70988
$ws.Range("J133").Value = 70988
$ws.Range("L133").Value = 70988
$ws.Range("N133").Value = -81108

$ws.Range("H136").Value = 3629.262
$ws.Range("I136").Value = 1589.5
$ws.Range("J136").Value = 5159.0835
$ws.Range("K136").Value = 4768.5
$ws.Range("L136").Value = 15477.2505
$ws.Range("M136").Value = -2218.5
$ws.Range("N136").Value = -20577.2505
